$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving/forcing Text cell type
# (matches the source workbook where Price/Volume columns are stored as inline strings)
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "70.845.35"
Set-TextCell $ws.Range("E2") "  -0.03%  "

# Row 3
Set-TextCell $ws.Range("D3") "3.529.36"
Set-TextCell $ws.Range("E3") "  -1.12%  "

# Row 4
Set-TextCell $ws.Range("D4") "0.999"
Set-TextCell $ws.Range("E4") "  -0.09%  "

# Row 5
Set-TextCell $ws.Range("D5") "613.06"
Set-TextCell $ws.Range("E5") "  -0.16%  "

# Row 6
Set-TextCell $ws.Range("D6") "173.57"
Set-TextCell $ws.Range("E6") "  +0.59%  "

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell $ws.Range("D7") "0.611"
Set-TextCell $ws.Range("E7") "  -1.42%  "

# Row 8
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextCell $ws.Range("D8") "3.523.59"
Set-TextCell $ws.Range("E8") "  -1.16%  "

# Row 9
Set-TextCell $ws.Range("E9") "  -0.07%  "

# Row 10
Set-TextCell $ws.Range("E10") "  -0.38%  "

# Row 11
Set-TextCell $ws.Range("D11") "7.32"
Set-TextCell $ws.Range("E11") "  +0.13%  "

# Row 13
Set-TextCell $ws.Range("D13") "46.59"
Set-TextCell $ws.Range("E13") "  -0.47%  "

# Row 14
Set-TextCell $ws.Range("E14") "  -0.46%  "

# Row 15
Set-TextCell $ws.Range("D15") "4.101.56"
Set-TextCell $ws.Range("E15") "  -0.93%  "

# Row 16
Set-TextCell $ws.Range("D16") "8.43"
Set-TextCell $ws.Range("E16") "  +0.19%  "

# Row 17
Set-TextCell $ws.Range("D17") "614.73"
Set-TextCell $ws.Range("E17") "  -0.69%  "

# Row 18
Set-TextCell $ws.Range("D18") "3.535.56"
Set-TextCell $ws.Range("E18") "  -0.83%  "

# Row 19
Set-TextCell $ws.Range("D19") "70.859.85"
Set-TextCell $ws.Range("E19") "  -0.14%  "

# Row 20
Set-TextCell $ws.Range("E20") "  +1.38%  "

# Row 21
Set-TextCell $ws.Range("D21") "17.77"
Set-TextCell $ws.Range("E21") "  +2.06%  "

# Row 22
Set-TextCell $ws.Range("D22") "0.885"
Set-TextCell $ws.Range("E22") "  +0.19%  "

# Row 23
Set-TextCell $ws.Range("D23") "8.99"
Set-TextCell $ws.Range("E23") "  -5.36%  "

# Row 24
Set-TextCell $ws.Range("D24") "15.75"
Set-TextCell $ws.Range("E24") "  -0.52%  "

# Row 25
Set-TextCell $ws.Range("D25") "98.07"
Set-TextCell $ws.Range("E25") "  +1.09%  "

# Row 26
Set-TextCell $ws.Range("E26") "  -1.54%  "

# Row 27
Set-TextCell $ws.Range("D27") "0.999"
Set-TextCell $ws.Range("E27") "  -0.07%  "

# Row 28
Set-TextCell $ws.Range("D28") "2.60"
Set-TextCell $ws.Range("E28") "  -0.61%  "

# Row 29
Set-TextCell $ws.Range("D29") "33.82"
Set-TextCell $ws.Range("E29") "  +0.32%  "

# Row 30
Set-TextCell $ws.Range("D30") "9.14"
Set-TextCell $ws.Range("E30") "  +0.24%  "

# Row 31
Set-TextCell $ws.Range("E31") "  -1.33%  "

# Row 32
Set-TextCell $ws.Range("D32") "8.16"
Set-TextCell $ws.Range("E32") "  -4.32%  "

# Row 33
Set-TextCell $ws.Range("E33") "  -0.58%  "

# Row 34
Set-TextCell $ws.Range("D34") "6.86"
Set-TextCell $ws.Range("E34") "  -1.57%  "

# Row 35
Set-TextCell $ws.Range("D35") "609.96"
Set-TextCell $ws.Range("E35") "  +6.54%  "

# Row 36
Set-TextCell $ws.Range("E36") "  -0.74%  "

# Row 37
Set-TextCell $ws.Range("D37") "10.85"
Set-TextCell $ws.Range("E37") "  -0.22%  "

# Row 38
Set-TextCell $ws.Range("D38") "3.54"
Set-TextCell $ws.Range("E38") "  -2.51%  "

# Row 39
Set-TextCell $ws.Range("D39") "0.0472"
Set-TextCell $ws.Range("E39") "  -0.30%  "

# Row 40
Set-TextCell $ws.Range("D40") "57.00"
Set-TextCell $ws.Range("E40") "  -1.03%  "

# Row 41
Set-TextCell $ws.Range("E41") "  +0.15%  "

# Row 42
Set-TextCell $ws.Range("E42") "  +1.24%  "

# Row 43
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws.Range("D43") "0.0₃0742"
Set-TextCell $ws.Range("E43") "  +5.23%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Range("D44") "3.366.76"
Set-TextCell $ws.Range("E44") "  -0.16%  "

# Row 45
Set-TextCell $ws.Range("D45") "0.313"
Set-TextCell $ws.Range("E45") "  -2.27%  "

# Row 46
Set-TextCell $ws.Range("E46") "  -2.01%  "

# Row 47
Set-TextCell $ws.Range("D47") "32.28"
Set-TextCell $ws.Range("E47") "  -2.36%  "

# Row 48
Set-TextCell $ws.Range("E48") "  -2.10%  "

# Row 49
Set-TextCell $ws.Range("E49") "  +0.42%  "

# Row 50
Set-TextCell $ws.Range("D50") "134.01"
Set-TextCell $ws.Range("E50") "  +0.18%  "
